# Generate Report for Handback
# - renames the UUID-named handback file (82a0a388... -> 2138835f...)
# - updates its hash/timestamps
# - appends a brand-new handback file row (49dce285...) to every sheet

$wb = $excel.ActiveWorkbook

$oldUuid = "82a0a388-0d8e-455e-9b26-ad91f2b58d3c"
$newUuid = "2138835f-8b04-4a55-8813-9128449d3f6c"
$addUuid = "49dce285-79b1-46d0-93eb-392e5b27552e"

$newHash = "98248e03a727e199df8d4695b38b4c67be377d70"
$addHash = "858ab1fd6d96947dbc11d5c84056b0ac6ca73c52"

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# row 2: rename the UUID throughout (value + hyperlink)
$wsOv.Range("A2").Value = "$newUuid.md"
$wsOv.Range("B2").Hyperlinks.Delete()
$wsOv.Range("B2").Value = "e2e\$newUuid.md"
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/812df44cdc10d13b2cfd11095ab98148ddae40f6/e2e/$newUuid.md", "", "", "e2e\$newUuid.md") | Out-Null
$wsOv.Range("G2").Value = "2016-09-03 15:06:31"

# row 3: brand-new file
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null
$wsOv.Range("A3").Value = "$addUuid.md"
$wsOv.Range("B3").Value = "e2e\$addUuid.md"
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/812df44cdc10d13b2cfd11095ab98148ddae40f6/e2e/$addUuid.md", "", "", "e2e\$addUuid.md") | Out-Null
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = $handedBack
$wsOv.Range("F3").Value = $handedBack
$wsOv.Range("G3").Value = "2016-09-03 15:06:31"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# row 2: rename UUID + refresh hash / timestamps for the existing handback
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("A2").Value = "$newUuid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/812df44cdc10d13b2cfd11095ab98148ddae40f6/e2e/$newUuid.md", "", "", "$newUuid.md") | Out-Null

$wsZh.Range("G2").Value = "$newUuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-03 15:06:26"

$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Range("I2").Value = "$newUuid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/46bcec4b89ac8f767c2fe91f7ff561d6a2823a65/e2e/$newUuid.md", "", "", "$newUuid.md") | Out-Null

$wsZh.Range("J2").Value = "$newUuid.$newHash.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-03 15:06:44"

# row 3: brand-new handback file
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = "$addUuid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/812df44cdc10d13b2cfd11095ab98148ddae40f6/e2e/$addUuid.md", "", "", "$addUuid.md") | Out-Null
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $handedBack
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$addUuid.$addHash.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-03 15:06:26"
$wsZh.Range("I3").Value = "$addUuid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/46bcec4b89ac8f767c2fe91f7ff561d6a2823a65/e2e/$addUuid.md", "", "", "$addUuid.md") | Out-Null
$wsZh.Range("J3").Value = "$addUuid.$addHash.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-03 15:06:44"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# row 2: rename UUID + refresh hash / timestamps for the existing handback
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("A2").Value = "$newUuid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/812df44cdc10d13b2cfd11095ab98148ddae40f6/e2e/$newUuid.md", "", "", "$newUuid.md") | Out-Null

$wsDe.Range("G2").Value = "$newUuid.$newHash.de-de.xlf"

$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Range("I2").Value = "$newUuid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/719c5287984539af34144bec95441fc8f3ee6aaf/e2e/$newUuid.md", "", "", "$newUuid.md") | Out-Null

$wsDe.Range("J2").Value = "$newUuid.$newHash.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 15:06:51"

# row 3: brand-new handback file
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = "$addUuid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/812df44cdc10d13b2cfd11095ab98148ddae40f6/e2e/$addUuid.md", "", "", "$addUuid.md") | Out-Null
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $handedBack
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$addUuid.$addHash.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-03 15:06:31"
$wsDe.Range("I3").Value = "$addUuid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/719c5287984539af34144bec95441fc8f3ee6aaf/e2e/$addUuid.md", "", "", "$addUuid.md") | Out-Null
$wsDe.Range("J3").Value = "$addUuid.$addHash.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-03 15:06:51"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

Write-Output "Handback report generated"
